{"js": "// Update the division problems (e.g. \"543\u00f75=\") inside the table cells.\n// The mapping below lists the old text followed by its replacement, in\n// the same order the paragraphs occur in the document body (this also\n// matches the order of the hunks in the source diff). Walking the\n// paragraphs once, left-to-right/top-to-bottom, and consuming the\n// mapping list in lock-step guarantees each paragraph is matched against\n// the correct pending replacement even though some \"new\" values are\n// identical to \"old\" values used later in the document\n// (e.g. \"363\u00f73=\" -> \"752\u00f76=\" and, further down, \"752\u00f76=\" -> \"478\u00f73=\").\nconst replacements = [\n  [\"543\u00f75=\", \"899\u00f72=\"],\n  [\"434\u00f76=\", \"498\u00f77=\"],\n  [\"122\u00f76=\", \"353\u00f75=\"],\n  [\"997\u00f77=\", \"563\u00f78=\"],\n  [\"913\u00f77=\", \"176\u00f79=\"],\n  [\"275\u00f77=\", \"955\u00f76=\"],\n  [\"720\u00f75=\", \"889\u00f79=\"],\n  [\"468\u00f72=\", \"372\u00f78=\"],\n  [\"849\u00f78=\", \"330\u00f76=\"],\n  [\"503\u00f74=\", \"921\u00f73=\"],\n  [\"525\u00f72=\", \"362\u00f76=\"],\n  [\"363\u00f73=\", \"752\u00f76=\"],\n  [\"123\u00f72=\", \"303\u00f73=\"],\n  [\"993\u00f75=\", \"184\u00f75=\"],\n  [\"712\u00f72=\", \"694\u00f77=\"],\n  [\"317\u00f77=\", \"164\u00f79=\"],\n  [\"367\u00f79=\", \"727\u00f79=\"],\n  [\"533\u00f75=\", \"457\u00f76=\"],\n  [\"654\u00f77=\", \"430\u00f75=\"],\n  [\"899\u00f77=\", \"381\u00f75=\"],\n  [\"252\u00f76=\", \"768\u00f79=\"],\n  [\"217\u00f78=\", \"756\u00f79=\"],\n  [\"347\u00f78=\", \"608\u00f74=\"],\n  [\"752\u00f76=\", \"478\u00f73=\"],\n  [\"247\u00f77=\", \"219\u00f75=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet next = 0;\nfor (let i = 0; i < paragraphs.items.length && next < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = replacements[next];\n  if (para.text === oldText) {\n    para.insertText(newText, \"Replace\");\n    next++;\n  }\n}\n\nawait context.sync();\n\nif (next !== replacements.length) {\n  throw new Error(\n    `Only matched ${next} of ${replacements.length} expected paragraphs`\n  );\n}\n", "ps1": "# Update the division problems (e.g. \"543\u00f75=\") inside the table cells.\n#\n# The list below pairs each old value with its replacement, in the same\n# order the paragraphs occur in the document (this also matches the\n# order of the hunks in the source diff). We walk every paragraph once,\n# top-to-bottom, and consume the pending-replacement list in lock-step so\n# each paragraph is matched against the correct replacement even though\n# some \"new\" values are identical to \"old\" values used later in the\n# document (e.g. \"363\u00f73=\" -> \"752\u00f76=\" and, further down,\n# \"752\u00f76=\" -> \"478\u00f73=\"). A naive global Find/Replace-All keyed only on\n# text would double-replace that shared value, so we replace one\n# paragraph at a time instead.\n$replacements = @(\n    @(\"543\u00f75=\", \"899\u00f72=\"),\n    @(\"434\u00f76=\", \"498\u00f77=\"),\n    @(\"122\u00f76=\", \"353\u00f75=\"),\n    @(\"997\u00f77=\", \"563\u00f78=\"),\n    @(\"913\u00f77=\", \"176\u00f79=\"),\n    @(\"275\u00f77=\", \"955\u00f76=\"),\n    @(\"720\u00f75=\", \"889\u00f79=\"),\n    @(\"468\u00f72=\", \"372\u00f78=\"),\n    @(\"849\u00f78=\", \"330\u00f76=\"),\n    @(\"503\u00f74=\", \"921\u00f73=\"),\n    @(\"525\u00f72=\", \"362\u00f76=\"),\n    @(\"363\u00f73=\", \"752\u00f76=\"),\n    @(\"123\u00f72=\", \"303\u00f73=\"),\n    @(\"993\u00f75=\", \"184\u00f75=\"),\n    @(\"712\u00f72=\", \"694\u00f77=\"),\n    @(\"317\u00f77=\", \"164\u00f79=\"),\n    @(\"367\u00f79=\", \"727\u00f79=\"),\n    @(\"533\u00f75=\", \"457\u00f76=\"),\n    @(\"654\u00f77=\", \"430\u00f75=\"),\n    @(\"899\u00f77=\", \"381\u00f75=\"),\n    @(\"252\u00f76=\", \"768\u00f79=\"),\n    @(\"217\u00f78=\", \"756\u00f79=\"),\n    @(\"347\u00f78=\", \"608\u00f74=\"),\n    @(\"752\u00f76=\", \"478\u00f73=\"),\n    @(\"247\u00f77=\", \"219\u00f75=\")\n)\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$next = 0\n\nfor ($i = 1; $i -le $paras.Count -and $next -lt $replacements.Count; $i++) {\n    $p = $paras.Item($i)\n    $range = $p.Range\n    $text = $range.Text.TrimEnd([char]13, [char]7)\n    $pair = $replacements[$next]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    if ($text -eq $oldText) {\n        $charRange = $d.Range($range.Start, $range.Start + $text.Length)\n        $charRange.Text = $newText\n        $next = $next + 1\n    }\n}\n\nif ($next -ne $replacements.Count) {\n    throw \"Only matched $next of $($replacements.Count) expected paragraphs\"\n}\n"}
